$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells P1, Q1 (continuing the 0..15 sequence), matching the
# same bold/centered/top-aligned/bordered formatting already used by the
# other header cells (B1:O1) so it resolves to the identical cell style.
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("P1:Q1").Font.Bold = $true
$ws.Range("P1:Q1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("P1:Q1").VerticalAlignment = -4160     # xlTop
$ws.Range("P1:Q1").Borders.LineStyle = 1         # xlContinuous

# Update I, K, M, O columns for rows 2-25 (swap the 1/2 pattern) and add
# the new P, Q columns (both filled with 2) for each of those rows
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I: was 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: was 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: was 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: was 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P: new column = 2
    $ws.Cells.Item($r, 17).Value = 2  # Q: new column = 2
}
